# Update header labels on each sheet so the first row reads better as a
# Power BI auto-detected header when the table is loaded.

$wb = $excel.ActiveWorkbook

# Sheets 1-3: "Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio (MWMed)",
# "Atendimento a Ponta(MW)" -> prefix "Ano " on B1:E1 (2015/2030/2040/2050)
foreach ($name in @("Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio (MWMed)", "Atendimento a Ponta(MW)")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("B1").Value = "Ano 2015"
    $ws.Range("C1").Value = "Ano 2030"
    $ws.Range("D1").Value = "Ano 2040"
    $ws.Range("E1").Value = "Ano 2050"
}

# Sheet 4: "Potencia Incremental - SIN(MW)" -> prefix "Intervalo " on B1:E1
$ws = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$ws.Range("B1").Value = "Intervalo 2015"
$ws.Range("C1").Value = "Intervalo 2015-2030"
$ws.Range("D1").Value = "Intervalo 2031-2040"
$ws.Range("E1").Value = "Intervalo 2041-2050"

# Sheet 5: "Emissoes Totais (MtCO2eq)" -> prefix "Ano " on B1:E1
$ws = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")
$ws.Range("B1").Value = "Ano 2015"
$ws.Range("C1").Value = "Ano 2030"
$ws.Range("D1").Value = "Ano 2040"
$ws.Range("E1").Value = "Ano 2050"

# Sheet 6: "Custo Total (bilhões de R$)" -> prefix "Ano " on B1 only
$ws = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws.Range("B1").Value = "Ano 2015"
